$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 24
$ws.Range("F3").Value = 2753
$ws.Range("F5").Value = 19986
$ws.Range("F6").Value = 81
$ws.Range("F7").Value = 2319
$ws.Range("F8").Value = 757
$ws.Range("F10").Value = 452
$ws.Range("F11").Value = 702
$ws.Range("F12").Value = 251
$ws.Range("F13").Value = 254
$ws.Range("F15").Value = 382
$ws.Range("F16").Value = 84
$ws.Range("F17").Value = 275
$ws.Range("F19").Value = 215
$ws.Range("F20").Value = 20
$ws.Range("F22").Value = 102

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 92
$ws.Range("F16").Value = 101

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6019
$ws.Range("F3").Value = 655

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6019
$ws.Range("F3").Value = 655
$ws.Range("F6").Value = 24
$ws.Range("F8").Value = 2753
$ws.Range("F10").Value = 19986
$ws.Range("F13").Value = 81
$ws.Range("F16").Value = 2319
$ws.Range("F17").Value = 757
$ws.Range("F20").Value = 452
$ws.Range("F21").Value = 702
$ws.Range("F22").Value = 251
$ws.Range("F23").Value = 254
$ws.Range("F28").Value = 382
$ws.Range("F29").Value = 84
$ws.Range("F32").Value = 275
$ws.Range("F33").Value = 92
$ws.Range("F36").Value = 215
$ws.Range("F37").Value = 101
$ws.Range("F38").Value = 101
$ws.Range("F40").Value = 20
$ws.Range("F49").Value = 102

